# predatorPrey_SPEED_TESTS - "Implemented predators into MAIN sim"
#
# Turns the single "BEFORE/AFTER REFACTORING" comparison block into a
# "FIRST ROUND" / "SECOND ROUND" comparison (before/after adding predators),
# and adds a spacer row (matching the note-row style) above the new
# "SECOND ROUND" section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the existing section headers in place (same cells, new text).
$ws.Range("B2").Value = "BEFORE ADDING PREDATORS"
$ws.Range("B19").Value = "SECOND ROUND"

# New header above everything: "FIRST ROUND" (row 1 was previously unused).
$ws.Range("B1").Value = "FIRST ROUND"

# Push the old "AFTER REFACTORING" template block (old rows 20-30) down one
# row, and use the freed-up row 20 for the new "AFTER ADDING PREDATORS" label.
$ws.Rows("20:20").Insert()
$ws.Range("B20").Value = "AFTER ADDING PREDATORS"

# Add an empty spacer cell at C18 (between the two notes in row 16/17 and the
# new "SECOND ROUND" row 19), carrying the same italic-note style as C17.
$ws.Range("C17").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The view had scrolled down with F20 selected; reset scroll to the top-left
# and select B20 (the new "AFTER ADDING PREDATORS" label cell).
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B20").Select()
